$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A1 value
$ws.Range("A1").Value = 842.28

# Delete rows 2 and 3 entirely (shifting cells up / removing them)
$ws.Range("A2:B3").EntireRow.Delete()
